# Tasks.xlsx - add "Wer" (who) owner column (D) values to the task list,
# per commit: "added smarter product catalog added search added
# increment/decrement product quantity added some minor improvements".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> owner name for column D (header "Wer" already in D2).
$owners = @{
    9  = "Kai"
    10 = "Kai"
    11 = "Kai"
    12 = "Kai"
    14 = "Kai"
    15 = "Partymarti"
    16 = "Partymarti"
    17 = "Partymarti"
    18 = "Riedo"
    19 = "Riedo"
    20 = "Riedo"
    21 = "Riedo"
    22 = "Riedo"
    23 = "Riedo"
    24 = "Alle"
    25 = "Riedo"
    26 = "Kai"
    27 = "Kai"
    28 = "Alle"
    29 = "Riedo"
    30 = "Alle"
    31 = "Kai"
    32 = "Kai"
    33 = "Kai"
    34 = "Kai"
    36 = "Riedo"
}

foreach ($row in $owners.Keys) {
    $ws.Cells.Item($row, 4).Value = $owners[$row]
}

# Leave the cursor/selection where the author last left it.
$ws.Range("B33").Select()
